$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing rows 2..157 down to 3..158,
# and automatically duplicates the last row's worth of data at the new end,
# since hyperlink refs/targets stay anchored to their original row numbers).
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the latest price entry. The date
# columns (A, E) hold plain text like "10-01-2026" in this sheet, not real
# Excel dates, so a leading apostrophe forces text entry instead of letting
# Excel's General-format autoconversion turn them into date serials.
$ws.Cells.Item(2, 1).Value = "'10-01-2026"
$ws.Cells.Item(2, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(2, 3).Value = "IE07"
$ws.Cells.Item(2, 4).Value = 307.25
$ws.Cells.Item(2, 5).Value = "'01-01-2026"
$ws.Cells.Item(2, 6).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-01-2026.pdf"

# The row-insert pushed the old last row's data to row 158, but left it
# without its own hyperlink relationship (hyperlinks don't move with the
# row shift in this engine). Register the hyperlink for the new row 158.
$ws.Hyperlinks.Add($ws.Range("F158"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf")
